$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.265357613563538
$ws.Range("B1").Value = 2.468292713165283
$ws.Range("C1").Value = 4.623887538909912
$ws.Range("D1").Value = 2.009847164154053
$ws.Range("E1").Value = 1.146043658256531
